$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "64.384.73"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.38%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.505.96"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "586.62"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "135.86"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.09%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.505.90"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +0.03%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.488"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -0.05%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.16"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.14%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.376"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.65%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.100.14"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  -0.08%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.499.94"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.29%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "64.366.04"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "24.96"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -10.49%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "9.80"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.83%  "
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("E21").Value = "  -2.38%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "383.91"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("E23").Value = "  -1.84%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.642.34"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "74.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E28").Value = "  +3.54%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.57"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.54"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.40%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.30"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  -1.35%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.524.76"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +0.52%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "23.57"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("E40").Value = "  -1.89%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "163.64"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.46%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0786"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.46%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.808"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.84%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "25.99"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("E45").Value = "  +0.08%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "41.89"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E49").Value = "  -0.70%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.478.14"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("E51").Value = "  -2.05%  "
